$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('kanji; Chinese character', '漢字|かんじ'),
    @('textbook', '教科書|きょうかしょ'),
    @('page', 'ページ'),
    @('next', '次|つぎ'),
    @('money', 'お金|おかね'),
    @('baggage', '荷物|にもつ'),
    @('personal computer', 'パソコン'),
    @('shower', 'シャワー'),
    @('air conditioner', 'エアコン'),
    @('electricity; light', '電気|でんき'),
    @('window', '窓|まど'),
    @('train', '電車|でんしゃ'),
    @('country; place of origin', '国|くに'),
    @('this week', '今週|こんしゅう'),
    @('next week', '来週|らいしゅう'),
    @('next year', '来年|らいねん'),
    @('night', '夜|よる'),
    @('tough (situation)', '大変|たいへん（な）'),
    @('to play; to spend time pleasantly', '遊ぶ|あそぶ'),
    @('to hurry', '急ぐ|いそぐ'),
    @('to return (a thing) (person に thing を)', '返す|かえす'),
    @('to turn off; to erase (～を)', '消す|けす'),
    @('to die', '死ぬ|しぬ'),
    @('to sit down (seat に)', '座る|すわる'),
    @('to stand up', '立つ|たつ'),
    @('to smoke', 'たばこを吸う|たばこをすう'),
    @('to use (～を)', '使う|つかう'),
    @('to help (person/task を)', '手伝う|てつだう'),
    @('to enter (～に)', '入る|はいる'),
    @('to carry; to hold (～を)', '持つ|もつ'),
    @('(1) to be absent (from...) (～を); (2) to rest', '休む|やすむ'),
    @('to open (something) (～を)', '開ける|あける'),
    @('to close (something) (～を)', '閉める|しめる'),
    @('to teach; to instruct (person に thing を)', '教える|おしえる'),
    @('to forget; to leave behind (～を)', '忘れる|わすれる'),
    @('to get off (～を)', '降りる|おりる'),
    @('to borrow (person に thing を)', '借りる|かりる'),
    @('to take a shower', 'シャワーを浴びる|シャワーをあびる'),
    @('to turn on (～を)', 'つける'),
    @('to call (～に)', '電話する|でんわする'),
    @('to bring (a person) (～を)', '連れてくる|つれてくる'),
    @('to bring (a thing) (～を)', '持ってくる|もってくる'),
    @('later on', '後で|あとで'),
    @('right away', 'すぐ'),
    @('slowly; leisurely; unhurriedly', 'ゆっくり'),
    @('That would be fine.; That wouldn''t be necessary.', '結構です|けっこうです'),
    @('Really?', '本当ですか|ほんとうですか'),
    @('to ask for directions; to ask the way', '道を聞く|みちをきく'),
    @('go straight', 'まっすぐ行く|まっすぐいき'),
    @('turn left', '左に曲がる|ひだりにまがる'),
    @('turn right', '右に曲がる|みぎにまがる'),
    @('cross the street', '道を渡る|みちをわたる'),
    @('turn left at the second corner', '二つ目の角を左に曲がる|ふたつめのかどをひだりにまがる'),
    @('turn right at the first traffic light', '一つ目の信号を右に曲がる|ひとつめのしんごうをみぎにまがる'),
    @('left side of the street', '道の左側|みちのひだりがわ'),
    @('right side of the street', '道の右側|みちのみぎがわ'),
    @('north', '北|きた'),
    @('east', '東|ひがし'),
    @('south', '南|みなみ'),
    @('west', '西|にし'),
    @('east', '東|ひがし'),
    @('east exit', '東口|ひがしぐち'),
    @('Tokyo', '東京|とうきょう'),
    @('Kanto area', '関東|かんとう'),
    @('the East', '東洋|とうよう'),
    @('west', '西|にし'),
    @('west exit', '西口|にしぐち'),
    @('northwest', '北西|ほくせい'),
    @('Kansai area', '関西|かんさい'),
    @('the West', '西洋|せいよう'),
    @('south', '南|みなみ'),
    @('south exit', '南口|みなみぐち'),
    @('southeast', '南東|なんとう'),
    @('Antarctica', '南極|なんきょく'),
    @('Southeast Asia', '東南アジア|とうなんアジア'),
    @('north', '北|きた'),
    @('north exit', '北口|きたぐち'),
    @('Tohoku area', '東北|とうほく'),
    @('North Pole', '北極|ほっきょく'),
    @('Hokkaido', '北海道|ほっかいどう'),
    @('mouth', '口|くち'),
    @('population', '人口|じんこう'),
    @('entrance', '入り口／入口|いりぐち'),
    @('to exit', '出る|でる'),
    @('exit', '出口|でぐち'),
    @('to go out', '出かける|でかける'),
    @('to take something out', '出す|だす'),
    @('attendance', '出席|しゅっせき'),
    @('export', '輸出|ゆしゅつ'),
    @('right', '右|みぎ'),
    @('right turn', '右折|うせつ'),
    @('right and left', '左右|さゆう'),
    @('right hand', '右手|みぎて'),
    @('right side', '右側|みぎがわ'),
    @('left', '左|ひだり'),
    @('left turn', '左折|させつ'),
    @('left hand', '左手|ひだりて'),
    @('left-handed', '左利き|ひだりきき'),
    @('five minutes', '五分|ごふん'),
    @('ten minutes', '十分|じゅっぷん／じっぷん'),
    @('oneself', '自分|じぶん'),
    @('to divide', '分ける|わける'),
    @('teacher', '先生|せんせい'),
    @('last week', '先週|せんしゅう'),
    @('ahead', '先に|さきに'),
    @('last month', '先月|せんげつ'),
    @('senior member', '先輩|せんぱい'),
    @('student', '学生|がくせい'),
    @('to be born', '生まれる|うまれる'),
    @('once in a life time', '一生に一度|いっしょうにいちど'),
    @('college student', '大学生|だいがくせい'),
    @('big', '大きい|おおきい'),
    @('tough', '大変な|たいへんな'),
    @('adult', '大人|おとな'),
    @('embassy', '大使館|たいしかん'),
    @('university', '大学|だいがく'),
    @('student', '学生|がくせい'),
    @('school', '学校|がっこう'),
    @('to study', '学ぶ|まなぶ'),
    @('department; faculty', '学部|がくぶ'),
    @('foreign country', '外国|がいこく'),
    @('foreigner', '外国人|がいこくじん'),
    @('outside', '外|そと'),
    @('overseas', '海外|かいがい'),
    @('China', '中国|ちゅうごく'),
    @('country', '国|くに'),
    @('South Korea', '韓国|かんこく'),
    @('the Diet', '国会|こっかい')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
}
